$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(2).Delete()
